$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = '59.485.26'; Numeric = $false },
    @{ Cell = "E2"; Value = '  -0.97%  '; Numeric = $false },
    @{ Cell = "D3"; Value = '2.640.72'; Numeric = $false },
    @{ Cell = "E3"; Value = '  +0.84%  '; Numeric = $false },
    @{ Cell = "E4"; Value = '  +0.10%  '; Numeric = $false },
    @{ Cell = "D5"; Value = '517.68'; Numeric = $true },
    @{ Cell = "E5"; Value = '  -0.58%  '; Numeric = $false },
    @{ Cell = "D6"; Value = '147.06'; Numeric = $true },
    @{ Cell = "E6"; Value = '  -1.04%  '; Numeric = $false },
    @{ Cell = "E7"; Value = '  -0.39%  '; Numeric = $false },
    @{ Cell = "E8"; Value = '  +0.72%  '; Numeric = $false },
    @{ Cell = "D9"; Value = '2.673.72'; Numeric = $false },
    @{ Cell = "E9"; Value = '  +1.97%  '; Numeric = $false },
    @{ Cell = "D10"; Value = '6.47'; Numeric = $true },
    @{ Cell = "E10"; Value = '  +2.44%  '; Numeric = $false },
    @{ Cell = "E11"; Value = '  +1.34%  '; Numeric = $false },
    @{ Cell = "E12"; Value = '  +0.18%  '; Numeric = $false },
    @{ Cell = "E13"; Value = '  -1.58%  '; Numeric = $false },
    @{ Cell = "D14"; Value = '3.109.61'; Numeric = $false },
    @{ Cell = "E14"; Value = '  +1.04%  '; Numeric = $false },
    @{ Cell = "D15"; Value = '59.463.63'; Numeric = $false },
    @{ Cell = "E15"; Value = '  -1.01%  '; Numeric = $false },
    @{ Cell = "D16"; Value = '21.23'; Numeric = $true },
    @{ Cell = "E16"; Value = '  +0.38%  '; Numeric = $false },
    @{ Cell = "E17"; Value = '  +0.57%  '; Numeric = $false },
    @{ Cell = "D18"; Value = '2.662.56'; Numeric = $false },
    @{ Cell = "E18"; Value = '  +1.71%  '; Numeric = $false },
    @{ Cell = "E19"; Value = '  -0.08%  '; Numeric = $false },
    @{ Cell = "D20"; Value = '346.47'; Numeric = $true },
    @{ Cell = "E20"; Value = '  +1.78%  '; Numeric = $false },
    @{ Cell = "D21"; Value = '10.52'; Numeric = $true },
    @{ Cell = "E21"; Value = '  +0.91%  '; Numeric = $false },
    @{ Cell = "E22"; Value = '  +1.12%  '; Numeric = $false },
    @{ Cell = "D23"; Value = '1.00'; Numeric = $true },
    @{ Cell = "E23"; Value = '  +0.58%  '; Numeric = $false },
    @{ Cell = "D24"; Value = '61.61'; Numeric = $true },
    @{ Cell = "E24"; Value = '  +1.41%  '; Numeric = $false },
    @{ Cell = "D25"; Value = '0.424'; Numeric = $true },
    @{ Cell = "E25"; Value = '  +1.27%  '; Numeric = $false },
    @{ Cell = "D26"; Value = '2.765.56'; Numeric = $false },
    @{ Cell = "E26"; Value = '  +0.89%  '; Numeric = $false },
    @{ Cell = "D27"; Value = '0.993'; Numeric = $true },
    @{ Cell = "E27"; Value = '  -0.31%  '; Numeric = $false },
    @{ Cell = "E28"; Value = '  +0.41%  '; Numeric = $false },
    @{ Cell = "E29"; Value = '  +1.80%  '; Numeric = $false },
    @{ Cell = "E30"; Value = '  +2.39%  '; Numeric = $false },
    @{ Cell = "D32"; Value = '6.49'; Numeric = $true },
    @{ Cell = "E32"; Value = '  +8.84%  '; Numeric = $false },
    @{ Cell = "E33"; Value = '  +0.88%  '; Numeric = $false },
    @{ Cell = "E34"; Value = '  -0.12%  '; Numeric = $false },
    @{ Cell = "D35"; Value = '150.00'; Numeric = $true },
    @{ Cell = "E35"; Value = '  +0.10%  '; Numeric = $false },
    @{ Cell = "E36"; Value = '  +13.75%  '; Numeric = $false },
    @{ Cell = "E37"; Value = '  +3.23%  '; Numeric = $false },
    @{ Cell = "E38"; Value = '  +3.41%  '; Numeric = $false },
    @{ Cell = "D39"; Value = '0.870'; Numeric = $true },
    @{ Cell = "E39"; Value = '  +1.08%  '; Numeric = $false },
    @{ Cell = "D40"; Value = '36.76'; Numeric = $true },
    @{ Cell = "E40"; Value = '  +0.79%  '; Numeric = $false },
    @{ Cell = "E41"; Value = '  +3.13%  '; Numeric = $false },
    @{ Cell = "E42"; Value = '  +0.30%  '; Numeric = $false },
    @{ Cell = "D43"; Value = '289.59'; Numeric = $true },
    @{ Cell = "E43"; Value = '  +0.88%  '; Numeric = $false },
    @{ Cell = "D44"; Value = '0.619'; Numeric = $true },
    @{ Cell = "E44"; Value = '  -0.86%  '; Numeric = $false },
    @{ Cell = "D45"; Value = '0.0993'; Numeric = $true },
    @{ Cell = "E45"; Value = '  -0.91%  '; Numeric = $false },
    @{ Cell = "D46"; Value = '0.995'; Numeric = $true },
    @{ Cell = "E46"; Value = '  -0.38%  '; Numeric = $false },
    @{ Cell = "D47"; Value = '19.79'; Numeric = $true },
    @{ Cell = "E47"; Value = '  +1.95%  '; Numeric = $false },
    @{ Cell = "D48"; Value = '0.0544'; Numeric = $true },
    @{ Cell = "E48"; Value = '  -0.27%  '; Numeric = $false },
    @{ Cell = "E49"; Value = '  +1.01%  '; Numeric = $false },
    @{ Cell = "E50"; Value = '  +1.80%  '; Numeric = $false },
    @{ Cell = "E51"; Value = '  -1.32%  '; Numeric = $false }
)

foreach ($item in $changes) {
    $rng = $ws.Range($item.Cell)
    if ($item.Numeric) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $item.Value
    if ($item.Numeric) {
        $rng.Style = "Normal"
    }
}
